$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: город Курск
$ws.Range("A6").Value = 38701000
$ws.Range("B6").Value = "город Курск"
$ws.Range("C6").Value = 2022
$ws.Range("D6").Value = -1897
$ws.Range("E6").Value = 447387
$ws.Range("F6").Value = 0.24330836613491227
$ws.Range("G6").Value = 25219.131849999998
$ws.Range("H6").Value = 0.89471419598691959
$ws.Range("I6").Value = 0.027714260807757043
$ws.Range("J6").Value = 97.734390398134039
$ws.Range("K6").Value = 30.6
$ws.Range("L6").Value = 0.0012517127229892688
$ws.Range("M6").Value = 0.0013455911772134639
$ws.Range("N6").Value = 0.0016786361695802516
$ws.Range("O6").Value = 0.0099645273555109999
$ws.Range("P6").Value = 0.13614680355039371
$ws.Range("Q6").Value = 1.272602612503269
$ws.Range("R6").Value = 0.0006482083744051571
$ws.Range("S6").Value = 0.036592480335816642
$ws.Range("T6").Value = 230.88501684045355

# Row 7: Новохопёрский МР
$ws.Range("A7").Value = 20627000
$ws.Range("B7").Value = "Новохопёрский МР"
$ws.Range("C7").Value = 2012
$ws.Range("D7").Value = -47
$ws.Range("E7").Value = 40333
$ws.Range("F7").Value = 0.11754642600352069
$ws.Range("G7").Value = 19539.2
$ws.Range("H7").Value = 0.35174794833015149
$ws.Range("I7").Value = 0.010611657947586344
$ws.Range("J7").Value = 5.3734534996156995
$ws.Range("K7").Value = 27.2
$ws.Range("L7").Value = 0.0030991991669352638
$ws.Range("M7").Value = 0.0012892668534450699
$ws.Range("N7").Value = 0.010331490342895395
$ws.Range("O7").Value = 2.2278531227530807
$ws.Range("P7").Value = 0.93820940668931152
$ws.Range("Q7").Value = 82.169960082314731
$ws.Range("R7").Value = 0.00079339498673542752
$ws.Range("S7").Value = 0.023206803362011256
$ws.Range("T7").Value = 101.40950586368483

# Apply styles matching row 4 pattern: columns A-E use style index 4, F-T use style index 5
$ws.Range("A6:E6").Style = $ws.Range("A4:E4").Style
$ws.Range("F6:T6").Style = $ws.Range("F4:T4").Style
$ws.Range("A7:E7").Style = $ws.Range("A4:E4").Style
$ws.Range("F7:T7").Style = $ws.Range("F4:T4").Style
